$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cell = $t.Cell(1,1)
$r = $cell.Range
$target = $d.Range($r.Start, $r.End - 1)
if ($target.Text -ne "206×7=") { Write-Host "MISMATCH at Row 1 Col 1: expected '206×7=' got '$($target.Text)'" }
$target.Text = "735×9="

$cell = $t.Cell(1,2)
$r = $cell.Range
$target = $d.Range($r.Start, $r.End - 1)
if ($target.Text -ne "747×2=") { Write-Host "MISMATCH at Row 1 Col 2: expected '747×2=' got '$($target.Text)'" }
$target.Text = "987×4="

$cell = $t.Cell(1,3)
$r = $cell.Range
$target = $d.Range($r.Start, $r.End - 1)
if ($target.Text -ne "113×7=") { Write-Host "MISMATCH at Row 1 Col 3: expected '113×7=' got '$($target.Text)'" }
$target.Text = "435×5="

$cell = $t.Cell(1,4)
$r = $cell.Range
$target = $d.Range($r.Start, $r.End - 1)
if ($target.Text -ne "549×3=") { Write-Host "MISMATCH at Row 1 Col 4: expected '549×3=' got '$($target.Text)'" }
$target.Text = "156×4="

$cell = $t.Cell(1,5)
$r = $cell.Range
$target = $d.Range($r.Start, $r.End - 1)
if ($target.Text -ne "446×5=") { Write-Host "MISMATCH at Row 1 Col 5: expected '446×5=' got '$($target.Text)'" }
$target.Text = "408×8="

$cell = $t.Cell(5,1)
$r = $cell.Range
$target = $d.Range($r.Start, $r.End - 1)
if ($target.Text -ne "183×8=") { Write-Host "MISMATCH at Row 5 Col 1: expected '183×8=' got '$($target.Text)'" }
$target.Text = "289×8="

$cell = $t.Cell(5,2)
$r = $cell.Range
$target = $d.Range($r.Start, $r.End - 1)
if ($target.Text -ne "838×5=") { Write-Host "MISMATCH at Row 5 Col 2: expected '838×5=' got '$($target.Text)'" }
$target.Text = "361×8="

$cell = $t.Cell(5,3)
$r = $cell.Range
$target = $d.Range($r.Start, $r.End - 1)
if ($target.Text -ne "946×3=") { Write-Host "MISMATCH at Row 5 Col 3: expected '946×3=' got '$($target.Text)'" }
$target.Text = "936×4="

$cell = $t.Cell(5,4)
$r = $cell.Range
$target = $d.Range($r.Start, $r.End - 1)
if ($target.Text -ne "467×2=") { Write-Host "MISMATCH at Row 5 Col 4: expected '467×2=' got '$($target.Text)'" }
$target.Text = "750×3="

$cell = $t.Cell(5,5)
$r = $cell.Range
$target = $d.Range($r.Start, $r.End - 1)
if ($target.Text -ne "668×6=") { Write-Host "MISMATCH at Row 5 Col 5: expected '668×6=' got '$($target.Text)'" }
$target.Text = "465×2="

$cell = $t.Cell(10,1)
$r = $cell.Range
$target = $d.Range($r.Start, $r.End - 1)
if ($target.Text -ne "987×6=") { Write-Host "MISMATCH at Row 10 Col 1: expected '987×6=' got '$($target.Text)'" }
$target.Text = "342×4="

$cell = $t.Cell(10,2)
$r = $cell.Range
$target = $d.Range($r.Start, $r.End - 1)
if ($target.Text -ne "957×7=") { Write-Host "MISMATCH at Row 10 Col 2: expected '957×7=' got '$($target.Text)'" }
$target.Text = "385×4="

$cell = $t.Cell(10,3)
$r = $cell.Range
$target = $d.Range($r.Start, $r.End - 1)
if ($target.Text -ne "814×6=") { Write-Host "MISMATCH at Row 10 Col 3: expected '814×6=' got '$($target.Text)'" }
$target.Text = "840×6="

$cell = $t.Cell(10,4)
$r = $cell.Range
$target = $d.Range($r.Start, $r.End - 1)
if ($target.Text -ne "456×4=") { Write-Host "MISMATCH at Row 10 Col 4: expected '456×4=' got '$($target.Text)'" }
$target.Text = "460×2="

$cell = $t.Cell(10,5)
$r = $cell.Range
$target = $d.Range($r.Start, $r.End - 1)
if ($target.Text -ne "574×9=") { Write-Host "MISMATCH at Row 10 Col 5: expected '574×9=' got '$($target.Text)'" }
$target.Text = "204×2="

$cell = $t.Cell(15,1)
$r = $cell.Range
$target = $d.Range($r.Start, $r.End - 1)
if ($target.Text -ne "427×7=") { Write-Host "MISMATCH at Row 15 Col 1: expected '427×7=' got '$($target.Text)'" }
$target.Text = "279×4="

$cell = $t.Cell(15,2)
$r = $cell.Range
$target = $d.Range($r.Start, $r.End - 1)
if ($target.Text -ne "881×7=") { Write-Host "MISMATCH at Row 15 Col 2: expected '881×7=' got '$($target.Text)'" }
$target.Text = "121×9="

$cell = $t.Cell(15,3)
$r = $cell.Range
$target = $d.Range($r.Start, $r.End - 1)
if ($target.Text -ne "364×3=") { Write-Host "MISMATCH at Row 15 Col 3: expected '364×3=' got '$($target.Text)'" }
$target.Text = "538×3="

$cell = $t.Cell(15,4)
$r = $cell.Range
$target = $d.Range($r.Start, $r.End - 1)
if ($target.Text -ne "142×6=") { Write-Host "MISMATCH at Row 15 Col 4: expected '142×6=' got '$($target.Text)'" }
$target.Text = "422×7="

$cell = $t.Cell(15,5)
$r = $cell.Range
$target = $d.Range($r.Start, $r.End - 1)
if ($target.Text -ne "457×9=") { Write-Host "MISMATCH at Row 15 Col 5: expected '457×9=' got '$($target.Text)'" }
$target.Text = "754×6="

$cell = $t.Cell(20,1)
$r = $cell.Range
$target = $d.Range($r.Start, $r.End - 1)
if ($target.Text -ne "129×8=") { Write-Host "MISMATCH at Row 20 Col 1: expected '129×8=' got '$($target.Text)'" }
$target.Text = "160×8="

$cell = $t.Cell(20,2)
$r = $cell.Range
$target = $d.Range($r.Start, $r.End - 1)
if ($target.Text -ne "573×2=") { Write-Host "MISMATCH at Row 20 Col 2: expected '573×2=' got '$($target.Text)'" }
$target.Text = "802×7="

$cell = $t.Cell(20,3)
$r = $cell.Range
$target = $d.Range($r.Start, $r.End - 1)
if ($target.Text -ne "909×4=") { Write-Host "MISMATCH at Row 20 Col 3: expected '909×4=' got '$($target.Text)'" }
$target.Text = "652×9="

$cell = $t.Cell(20,4)
$r = $cell.Range
$target = $d.Range($r.Start, $r.End - 1)
if ($target.Text -ne "881×7=") { Write-Host "MISMATCH at Row 20 Col 4: expected '881×7=' got '$($target.Text)'" }
$target.Text = "899×2="

$cell = $t.Cell(20,5)
$r = $cell.Range
$target = $d.Range($r.Start, $r.End - 1)
if ($target.Text -ne "559×8=") { Write-Host "MISMATCH at Row 20 Col 5: expected '559×8=' got '$($target.Text)'" }
$target.Text = "675×9="
